$p = $ppt.ActivePresentation

# Slide 1 title: collapse "First" / " " / "slide" runs into a single run "First slide"
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Text = "F"
$tr1.Text = "First slide"

# Slide 2 (blank slide) notes page text: collapse the many word/space runs into one run
$s2 = $p.Slides.Item(2)
$notes = $s2.NotesPage
$trN = $notes.Shapes.Item(2).TextFrame.TextRange
$trN.Text = "S"
$trN.Text = "Some notes here: this first slide should use the Blank template"

# Slide 3 title: collapse "Third" / " " / "slide" runs into a single run "Third slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Text = "T"
$tr3.Text = "Third slide"
